$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 88181

# Row 3
$ws.Range("B3").Value = 90448

# Row 4
$ws.Range("A4").Value = 112370020
$ws.Range("Q4").Value = 469262
$ws.Range("R4").Value = 7039652
$ws.Range("AC4").Value = "ringhack äldre"

# Row 5
$ws.Range("A5").Value = 112370021
$ws.Range("Q5").Value = 469287
$ws.Range("R5").Value = 7039645
$ws.Range("AC5").Value = "ringhack"

# Row 6
$ws.Range("A6").Value = 112067971
$ws.Range("B6").Value = 90799
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 1968
$ws.Range("F6").Value = "Grantaggsvamp"
$ws.Range("G6").Value = "Bankera violascens"
$ws.Range("H6").Value = "(Alb. & Schwein. : Fr.) Pouzar"

# Row 7
$ws.Range("A7").Value = 112068136
$ws.Range("B7").Value = 88180
$ws.Range("D7").Value = "VU"
$ws.Range("E7").Value = 6276
$ws.Range("F7").Value = "Goliatmusseron"
$ws.Range("G7").Value = "Tricholoma matsutake"
$ws.Range("H7").Value = "(S.Ito & S.Imai) Singer"
$ws.Range("Q7").Value = 469497
$ws.Range("R7").Value = 7039592

# Row 8
$ws.Range("A8").Value = 112068010
$ws.Range("B8").Value = 88180
$ws.Range("Q8").Value = 469452
$ws.Range("R8").Value = 7039595

# Row 9
$ws.Range("A9").Value = 112067953
$ws.Range("B9").Value = 88180
$ws.Range("D9").Value = "VU"
$ws.Range("E9").Value = 6276
$ws.Range("F9").Value = "Goliatmusseron"
$ws.Range("G9").Value = "Tricholoma matsutake"
$ws.Range("H9").Value = "(S.Ito & S.Imai) Singer"

# Row 10
$ws.Range("A10").Value = 112068040
$ws.Range("B10").Value = 90830
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 2059
$ws.Range("F10").Value = "Skrovlig taggsvamp"
$ws.Range("G10").Value = "Hydnellum scabrosum"
$ws.Range("H10").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q10").Value = 469465
$ws.Range("R10").Value = 7039571

# Row 11
$ws.Range("B11").Value = 90830
